# Apply cryptos list update (values scraped on Sat Sep  2 21:06:20 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.897.88"
$ws.Range("E2").Value = "  -0.61%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.639.12"
$ws.Range("E3").Value = "  -0.37%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.67%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.98"
$ws.Range("E5").Value = "  -0.40%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5025"
$ws.Range("E6").Value = "  -0.37%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.004"
$ws.Range("E7").Value = "  -0.74%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2565"
$ws.Range("E8").Value = "  -0.75%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06393"
$ws.Range("E9").Value = "  -1.08%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.60"
$ws.Range("E10").Value = "  +0.30%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07799"
$ws.Range("E11").Value = "  +0.59%  "

# Row 12
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.275"
$ws.Range("E12").Value = "  +0.06%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.650.43"
$ws.Range("E13").Value = "  +0.35%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.862.43"
$ws.Range("E14").Value = "  -0.27%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5430"
$ws.Range("E15").Value = "  -0.71%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅7858"
$ws.Range("E16").Value = "  -1.24%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.65"
$ws.Range("E17").Value = "  +1.36%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.907.64"
$ws.Range("E18").Value = "  -0.47%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.004"
$ws.Range("E19").Value = "  -0.59%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "198.28"
$ws.Range("E20").Value = "  -3.14%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.386"
$ws.Range("E21").Value = "  +1.54%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.936"
$ws.Range("E22").Value = "  -0.81%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.972"
$ws.Range("E23").Value = "  -0.08%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.003"
$ws.Range("E24").Value = "  -0.91%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.873"
$ws.Range("E25").Value = "  -3.36%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "139.90"
$ws.Range("E26").Value = "  -1.50%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1140"
$ws.Range("E27").Value = "  -1.49%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.839"
$ws.Range("E28").Value = "  +1.22%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.67"
$ws.Range("E29").Value = "  -1.15%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.242"
$ws.Range("E30").Value = "  -0.25%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04889"
$ws.Range("E31").Value = "  -3.76%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.256"
$ws.Range("E32").Value = "  -0.59%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.192"
$ws.Range("E33").Value = "  -0.33%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.530"
$ws.Range("E34").Value = "  -1.15%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.368"
$ws.Range("E35").Value = "  +0.76%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.8912"
$ws.Range("E36").Value = "  -0.75%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.599"
$ws.Range("E37").Value = "  -1.00%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.135.65"
$ws.Range("E38").Value = "  -1.45%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5544"
$ws.Range("E39").Value = "  -2.16%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01557"
$ws.Range("E40").Value = "  -0.91%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.006"
$ws.Range("E41").Value = "  -0.67%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.686"
$ws.Range("E42").Value = "  +0.40%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8155"
$ws.Range("E43").Value = "  -0.51%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.49"
$ws.Range("E44").Value = "  -0.29%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₈124"
$ws.Range("E45").Value = "  +11.02%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.773.79"
$ws.Range("E46").Value = "  -0.26%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4526"
$ws.Range("E47").Value = "  -0.46%  "

# Row 48
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "55.32"
$ws.Range("E48").Value = "  +0.46%  "

# Row 49
$ws.Range("B49").Value = "Frax"
$ws.Range("C49").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.006"
$ws.Range("E49").Value = "  -0.18%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05083"
$ws.Range("E50").Value = "  +0.70%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.005"
$ws.Range("E51").Value = "  -0.38%  "

